# run 9 and 10 TA for E5 samples dmb 0715
# Adds the new CRM-accuracy reading captured on 2021-07-15 as row 46,
# extending the existing "% off" shared formula down into the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row
$ws.Range("A46").Value = 20210715
$ws.Range("B46").Value = 2230.4602923572402
$ws.Range("C46").Value = 2224.4699999999998
$ws.Range("E46").Value = 180
$ws.Range("F46").Value = "CRM OPENED 20210526"

# Extend the "% off" formula (D40:D45 previously) down through the new row.
# Re-entering the formula across the whole D40:D46 block lets D41:D45 keep
# their existing results while D46 picks up the same relative formula.
$ws.Range("D40:D46").Formula = "=100*(B40-C40)/C40"

# Move the view down to the newly added row and select the new Notes cell,
# matching where Excel leaves the cursor after entering the row.
$ws.Range("F46").Select() | Out-Null
